$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source cells are text (inlineStr) even where the text looks numeric (e.g. "572.98").
# Excel COM auto-coerces plain decimal-looking strings into real numbers on assignment,
# so for those specific cells we set Text format first to keep them as strings - this
# matches the original file where every Coin/Link/Price/Volume cell is stored as text.
# Cells whose new text is unambiguous (multi-dot prices, %-strings, plain text) do not
# need this and are written directly.

$ws.Range('D2').Value = '60.439.17'
$ws.Range('E2').Value = '  +1.57%  '
$ws.Range('D3').Value = '2.602.39'
$ws.Range('E3').Value = '  +0.43%  '
$ws.Range('E4').Value = '  +0.00%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '572.98'
$ws.Range('E5').Value = '  +2.08%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '142.44'
$ws.Range('E6').Value = '  -0.75%  '
$ws.Range('E7').Value = '  -0.16%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.601'
$ws.Range('E8').Value = '  +0.50%  '
$ws.Range('D9').Value = '2.626.66'
$ws.Range('E9').Value = '  +0.87%  '
$ws.Range('E10').Value = '  -2.81%  '
$ws.Range('E11').Value = '  +0.82%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.153'
$ws.Range('E12').Value = '  -4.05%  '
$ws.Range('E13').Value = '  +2.13%  '
$ws.Range('D14').Value = '3.068.32'
$ws.Range('E14').Value = '  +0.58%  '
$ws.Range('D15').Value = '60.457.85'
$ws.Range('E15').Value = '  +1.67%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '23.23'
$ws.Range('E16').Value = '  -0.77%  '
$ws.Range('E17').Value = '  +2.36%  '
$ws.Range('D18').Value = '2.614.86'
$ws.Range('E18').Value = '  +0.75%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '11.34'
$ws.Range('E19').Value = '  +8.63%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '4.65'
$ws.Range('E20').Value = '  +1.34%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '346.43'
$ws.Range('E21').Value = '  +2.27%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '6.98'
$ws.Range('E22').Value = '  +6.35%  '
$ws.Range('E23').Value = '  -0.47%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '0.529'
$ws.Range('E24').Value = '  +12.20%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '63.26'
$ws.Range('E25').Value = '  -0.96%  '
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('E27').Value = '  -1.51%  '
$ws.Range('E28').Value = '  +3.61%  '
$ws.Range('E29').Value = '  +0.78%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '1.85'
$ws.Range('E30').Value = '  +10.31%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '6.38'
$ws.Range('E31').Value = '  +2.89%  '
$ws.Range('E32').Value = '  -0.10%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '161.32'
$ws.Range('E33').Value = '  +1.90%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '19.50'
$ws.Range('E34').Value = '  +1.92%  '
$ws.Range('E35').Value = '  +3.81%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.982'
$ws.Range('E36').Value = '  +9.51%  '
$ws.Range('E37').Value = '  +3.98%  '
$ws.Range('E38').Value = '  +7.44%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '37.83'
$ws.Range('E39').Value = '  +1.00%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '3.85'
$ws.Range('E40').Value = '  +4.22%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.847'
$ws.Range('E41').Value = '  -3.25%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '294.61'
$ws.Range('E42').Value = '  +0.12%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '137.64'
$ws.Range('E43').Value = '  -1.21%  '
$ws.Range('E44').Value = '  -0.34%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.609'
$ws.Range('E45').Value = '  +1.86%  '
$ws.Range('E46').Value = '  +0.74%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '19.79'
$ws.Range('E47').Value = '  +3.04%  '
$ws.Range('E48').Value = '  +2.62%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '4.93'
$ws.Range('E49').Value = '  +8.65%  '
$ws.Range('E50').Value = '  +1.73%  '
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '10.72'
$ws.Range('E51').Value = '  +0.73%  '
